$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "edit2"
$ws.Range("B3").Value = "riya-morankar"
$ws.Range("C3").Value = "Merged"
$ws.Range("D3").Value = "comment!"
# Force text format on the Date column so the "YYYY-MM-DD" string is kept
# literally instead of being auto-parsed into a date serial number (matches
# the plain-text value already used by every other row in this column).
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2025-06-18"
$ws.Range("F3").Value = "N/A"
